$wb = $excel.ActiveWorkbook
$target = $wb.Worksheets.Item("securitygroups")
$ws = $wb.Worksheets.Add($target)
$ws.Name = "volumes"

$ws.Range("A1").Value = "*name"
$ws.Range("B1").Value = "*zone"
$ws.Range("C1").Value = "*profile"
$ws.Range("D1").Value = "iops"
$ws.Range("E1").Value = "capacity"
$ws.Range("F1").Value = "encryption_key"
$ws.Range("G1").Value = "resource_group"
$ws.Range("H1").Value = "resource_controller_url"
$ws.Range("I1").Value = "create_timeout"
$ws.Range("J1").Value = "delete_timeout"

$ws.Range("A2").Value = "volume1"
$ws.Range("B2").Value = "Dallas 1"
$ws.Range("C2").Value = "10 IOPS/GB"
$ws.Range("E2").Value = 100
$ws.Range("I2").Value = "60m"
$ws.Range("J2").Value = "60m"

$instances = $wb.Worksheets.Item("instances")
$instances.Range("M2").Value = "volume1"
